# Weekly crime data refresh: shift reporting week forward by one week
# and update the 010 Pct CompStat figures (volume/issue number, week-of
# dates, and the weekly/28-day/YTD/2yr crime counts + % change columns).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Title block: bump issue number & the reporting week's date range ---
$ws.Range("A8").Value = "Volume 31   Number  36"
$ws.Range("C9").Value = "Report Covering the Week  9/2/2024  Through  9/8/2024"

# --- Crime complaint table (rows 15-31): updated counts / %-change ---
$ws.Range("F15").Value = "0"
$ws.Range("G15").Value = "0"
$ws.Range("H15").Value = "***.*"
$ws.Range("L15").Value = 42.857142857142
$ws.Range("N15").Value = -44.444444444444
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = -50
$ws.Range("F16").Value = 8
$ws.Range("G16").Value = 14
$ws.Range("H16").Value = -42.857142857142
$ws.Range("I16").Value = 73
$ws.Range("J16").Value = 99
$ws.Range("K16").Value = -26.262626262626
$ws.Range("L16").Value = -29.807692307692
$ws.Range("M16").Value = 5.797101449275
$ws.Range("N16").Value = -82.983682983683
$ws.Range("D17").Value = "0"
$ws.Range("E17").Value = "***.*"
$ws.Range("G17").Value = 7
$ws.Range("H17").Value = 57.142857142857
$ws.Range("I17").Value = 85
$ws.Range("K17").Value = -24.107142857142
$ws.Range("L17").Value = -5.555555555555
$ws.Range("M17").Value = -7.608695652173
$ws.Range("N17").Value = -50.581395348837
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 5
$ws.Range("G18").Value = 15
$ws.Range("H18").Value = -66.666666666666
$ws.Range("I18").Value = 70
$ws.Range("J18").Value = 93
$ws.Range("K18").Value = -24.731182795698
$ws.Range("L18").Value = -45.3125
$ws.Range("M18").Value = 1.449275362318
$ws.Range("N18").Value = -82.885085574572
$ws.Range("C19").Value = 11
$ws.Range("D19").Value = 14
$ws.Range("E19").Value = -21.428571428571
$ws.Range("F19").Value = 59
$ws.Range("G19").Value = 59
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 474
$ws.Range("J19").Value = 543
$ws.Range("K19").Value = -12.707182320442
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = 13.397129186602
$ws.Range("N19").Value = -15.508021390374
$ws.Range("C20").Value = 3
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = -25
$ws.Range("F20").Value = 8
$ws.Range("G20").Value = 13
$ws.Range("H20").Value = -38.461538461538
$ws.Range("I20").Value = 35
$ws.Range("J20").Value = 56
$ws.Range("K20").Value = -37.5
$ws.Range("L20").Value = -41.666666666666
$ws.Range("M20").Value = 9.375
$ws.Range("N20").Value = -89.489489489489
$ws.Range("C21").Value = 20
$ws.Range("D21").Value = 22
$ws.Range("E21").Value = -9.090909090909
$ws.Range("F21").Value = 91
$ws.Range("G21").Value = 108
$ws.Range("H21").Value = -15.740740740740
$ws.Range("I21").Value = 747
$ws.Range("J21").Value = 907
$ws.Range("K21").Value = -17.640573318632
$ws.Range("L21").Value = -13.741339491916
$ws.Range("M21").Value = 8.733624454148
$ws.Range("N21").Value = -61.235080435910
$ws.Range("C22").Value = "0"
$ws.Range("D22").Value = "0"
$ws.Range("E22").Value = "***.*"
$ws.Range("G22").Value = 3
$ws.Range("H22").Value = -66.666666666666
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = -100
$ws.Range("F23").Value = 2
$ws.Range("H23").Value = -33.333333333333
$ws.Range("J23").Value = 34
$ws.Range("K23").Value = -11.764705882352
$ws.Range("L23").Value = -28.571428571428
$ws.Range("M23").Value = -9.090909090909
$ws.Range("D24").Value = 9
$ws.Range("E24").Value = 44.444444444444
$ws.Range("F24").Value = 62
$ws.Range("G24").Value = 83
$ws.Range("H24").Value = -25.301204819277
$ws.Range("I24").Value = 469
$ws.Range("J24").Value = 570
$ws.Range("K24").Value = -17.719298245614
$ws.Range("L24").Value = -7.858546168958
$ws.Range("M24").Value = -22.861842105263
$ws.Range("C25").Value = 6
$ws.Range("D25").Value = 5
$ws.Range("E25").Value = 20
$ws.Range("F25").Value = 31
$ws.Range("G25").Value = 48
$ws.Range("H25").Value = -35.416666666666
$ws.Range("I25").Value = 204
$ws.Range("J25").Value = 340
$ws.Range("K25").Value = -40
$ws.Range("L25").Value = -30.136986301369
$ws.Range("C26").Value = 7
$ws.Range("D26").Value = 5
$ws.Range("E26").Value = 40
$ws.Range("F26").Value = 30
$ws.Range("G26").Value = 25
$ws.Range("H26").Value = 20
$ws.Range("I26").Value = 221
$ws.Range("J26").Value = 233
$ws.Range("K26").Value = -5.150214592274
$ws.Range("L26").Value = -14.007782101167
$ws.Range("M26").Value = -4.741379310344
$ws.Range("F27").Value = "0"
$ws.Range("G27").Value = "0"
$ws.Range("H27").Value = "***.*"
$ws.Range("L27").Value = 58.333333333333
$ws.Range("F28").Value = 4
$ws.Range("G28").Value = 4
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 41
$ws.Range("K28").Value = -2.380952380952
$ws.Range("L28").Value = 2.5
$ws.Range("D31").Value = "0"
$ws.Range("E31").Value = "***.*"
$ws.Range("L31").Value = -64.285714285714
